$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LUA CODING")

# Mark "Auto adjust Rep boxes when rank changes." and
# "Code Reputation to roll on a click button." as 100% done
# (commit: "Updated LUA to roll Reputation dice.")
$ws.Range("B15").Value = 1
$ws.Range("B16").Value = 1
$ws.Range("B17").Value = 1

# Update the view state to reflect scrolling down to the edited rows
$ws.Activate()
$ws.Range("B18").Select()
